# Auto-update draw results: append the 2025-12-24 Pick 3 draw as a new row
# at the bottom of the Results sheet (row 99), mirroring the existing
# table layout (Date, Game, Phase, Result, InsertedAt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and C hold values that look like a date / a pure number
# ("2025-12-24" and "251224"); force them to text first so Excel keeps
# them as literal strings instead of silently converting them to a
# date serial / numeric value (matching how every other row in the
# table is stored as text).
$ws.Range("A99").NumberFormat = "@"
$ws.Range("C99").NumberFormat = "@"

$ws.Range("A99").Value = "2025-12-24"
$ws.Range("B99").Value = "Pick 3"
$ws.Range("C99").Value = "251224"
$ws.Range("D99").Value = "0-8-1"
$ws.Range("E99").Value = "2025-12-24T21:39:33.813+04:00"
